# Update calculated price/profit columns (H-N) for specific Leve rows
# across each job sheet, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1953.3673
$ws.Cells.Item(132, 9).Value = 1026.5714
$ws.Cells.Item(132, 11).Value = 3079.7142
$ws.Cells.Item(132, 13).Value = -549.7142000000003
$ws.Cells.Item(138, 8).Value = 1989.2979
$ws.Cells.Item(138, 9).Value = 915.96295
$ws.Cells.Item(138, 10).Value = 3438.3
$ws.Cells.Item(138, 11).Value = 2747.88885
$ws.Cells.Item(138, 12).Value = 10314.9
$ws.Cells.Item(138, 13).Value = 2392.11115
$ws.Cells.Item(138, 14).Value = -20594.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 2500
$ws.Cells.Item(17, 9).Value = 2000
$ws.Cells.Item(17, 10).Value = 3000
$ws.Cells.Item(17, 11).Value = 2000
$ws.Cells.Item(17, 12).Value = 3000
$ws.Cells.Item(17, 13).Value = -1827
$ws.Cells.Item(17, 14).Value = -3346
$ws.Cells.Item(32, 8).Value = 1493260.1
$ws.Cells.Item(32, 9).Value = 1669671.9
$ws.Cells.Item(32, 11).Value = 1669671.9
$ws.Cells.Item(32, 13).Value = -1669384.9
$ws.Cells.Item(34, 8).Value = 11000
$ws.Cells.Item(34, 9).Value = 7000
$ws.Cells.Item(34, 11).Value = 7000
$ws.Cells.Item(34, 13).Value = -6729
$ws.Cells.Item(45, 8).Value = 1409.3636
$ws.Cells.Item(45, 9).Value = 1189.2222
$ws.Cells.Item(45, 11).Value = 1189.2222
$ws.Cells.Item(45, 13).Value = -812.2221999999999
$ws.Cells.Item(61, 8).Value = 387069.53
$ws.Cells.Item(61, 9).Value = 304814.06
$ws.Cells.Item(61, 10).Value = 529934.3
$ws.Cells.Item(61, 11).Value = 304814.06
$ws.Cells.Item(61, 12).Value = 529934.3
$ws.Cells.Item(61, 13).Value = -304602.06
$ws.Cells.Item(61, 14).Value = -530358.3
$ws.Cells.Item(74, 8).Value = 213195.94
$ws.Cells.Item(74, 9).Value = 271237.16
$ws.Cells.Item(74, 10).Value = 70027.60000000001
$ws.Cells.Item(74, 11).Value = 271237.16
$ws.Cells.Item(74, 12).Value = 70027.60000000001
$ws.Cells.Item(74, 13).Value = -270363.16
$ws.Cells.Item(74, 14).Value = -71775.60000000001
$ws.Cells.Item(77, 8).Value = 213195.94
$ws.Cells.Item(77, 9).Value = 271237.16
$ws.Cells.Item(77, 10).Value = 70027.60000000001
$ws.Cells.Item(77, 11).Value = 1356185.8
$ws.Cells.Item(77, 12).Value = 350138
$ws.Cells.Item(77, 13).Value = -1351817.8
$ws.Cells.Item(77, 14).Value = -358874
$ws.Cells.Item(136, 8).Value = 387069.53
$ws.Cells.Item(136, 9).Value = 304814.06
$ws.Cells.Item(136, 10).Value = 529934.3
$ws.Cells.Item(136, 11).Value = 914442.1799999999
$ws.Cells.Item(136, 12).Value = 1589802.9
$ws.Cells.Item(136, 13).Value = -911892.1799999999
$ws.Cells.Item(136, 14).Value = -1594902.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2735.244
$ws.Cells.Item(86, 9).Value = 2913.7878
$ws.Cells.Item(86, 11).Value = 2913.7878
$ws.Cells.Item(86, 13).Value = -1790.7878
$ws.Cells.Item(89, 8).Value = 2735.244
$ws.Cells.Item(89, 9).Value = 2913.7878
$ws.Cells.Item(89, 11).Value = 14568.939
$ws.Cells.Item(89, 13).Value = -8952.939
$ws.Cells.Item(105, 8).Value = 2858880.5
$ws.Cells.Item(105, 9).Value = 1631.5385
$ws.Cells.Item(105, 10).Value = 11113156
$ws.Cells.Item(105, 11).Value = 1631.5385
$ws.Cells.Item(105, 12).Value = 11113156
$ws.Cells.Item(105, 13).Value = 115.4614999999999
$ws.Cells.Item(105, 14).Value = -11116650
$ws.Cells.Item(107, 8).Value = 1925.8948
$ws.Cells.Item(107, 9).Value = 1819.2142
$ws.Cells.Item(107, 10).Value = 2224.6
$ws.Cells.Item(107, 11).Value = 1819.2142
$ws.Cells.Item(107, 12).Value = 2224.6
$ws.Cells.Item(107, 13).Value = 100.7858000000001
$ws.Cells.Item(107, 14).Value = -6064.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 860
$ws.Cells.Item(16, 9).Value = 670
$ws.Cells.Item(16, 10).Value = 1525
$ws.Cells.Item(16, 11).Value = 670
$ws.Cells.Item(16, 12).Value = 1525
$ws.Cells.Item(16, 13).Value = -383
$ws.Cells.Item(16, 14).Value = -2099
$ws.Cells.Item(99, 8).Value = 85026
$ws.Cells.Item(99, 10).Value = 2114.2856
$ws.Cells.Item(99, 12).Value = 2114.2856
$ws.Cells.Item(99, 14).Value = -5110.2856
$ws.Cells.Item(113, 8).Value = 860
$ws.Cells.Item(113, 9).Value = 670
$ws.Cells.Item(113, 10).Value = 1525
$ws.Cells.Item(113, 11).Value = 670
$ws.Cells.Item(113, 12).Value = 1525
$ws.Cells.Item(113, 13).Value = 1500
$ws.Cells.Item(113, 14).Value = -5865
$ws.Cells.Item(126, 8).Value = 85026
$ws.Cells.Item(126, 10).Value = 2114.2856
$ws.Cells.Item(126, 12).Value = 6342.8568
$ws.Cells.Item(126, 14).Value = -11282.8568
$ws.Cells.Item(132, 8).Value = 2239.5642
$ws.Cells.Item(132, 9).Value = 1389.8948
$ws.Cells.Item(132, 11).Value = 4169.6844
$ws.Cells.Item(132, 13).Value = -1639.6844
$ws.Cells.Item(134, 8).Value = 1516.8
$ws.Cells.Item(134, 9).Value = 654.9474
$ws.Cells.Item(134, 11).Value = 1964.8422
$ws.Cells.Item(134, 13).Value = 570.1578

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1335.7826
$ws.Cells.Item(5, 9).Value = 510.8889
$ws.Cells.Item(5, 10).Value = 1866.0714
$ws.Cells.Item(5, 11).Value = 1532.6667
$ws.Cells.Item(5, 12).Value = 5598.2142
$ws.Cells.Item(5, 13).Value = -1420.6667
$ws.Cells.Item(5, 14).Value = -5822.2142
$ws.Cells.Item(122, 8).Value = 25000584
$ws.Cells.Item(122, 9).Value = 33333798
$ws.Cells.Item(122, 10).Value = 939.6
$ws.Cells.Item(122, 11).Value = 300004182
$ws.Cells.Item(122, 12).Value = 8456.4
$ws.Cells.Item(122, 13).Value = -300001732
$ws.Cells.Item(122, 14).Value = -13356.4
$ws.Cells.Item(123, 8).Value = 8523.333000000001
$ws.Cells.Item(123, 9).Value = 13213.333
$ws.Cells.Item(123, 11).Value = 39639.999
$ws.Cells.Item(123, 13).Value = -37189.999
$ws.Cells.Item(131, 8).Value = 1166.4769
$ws.Cells.Item(131, 10).Value = 1254.1552
$ws.Cells.Item(131, 12).Value = 3762.4656
$ws.Cells.Item(131, 14).Value = -13842.4656
$ws.Cells.Item(132, 8).Value = 5698.1816
$ws.Cells.Item(132, 9).Value = 4133.091
$ws.Cells.Item(132, 10).Value = 7263.273
$ws.Cells.Item(132, 11).Value = 37197.819
$ws.Cells.Item(132, 12).Value = 65369.457
$ws.Cells.Item(132, 13).Value = -34667.819
$ws.Cells.Item(132, 14).Value = -70429.45699999999
$ws.Cells.Item(135, 8).Value = 1335.7826
$ws.Cells.Item(135, 9).Value = 510.8889
$ws.Cells.Item(135, 10).Value = 1866.0714
$ws.Cells.Item(135, 11).Value = 4598.0001
$ws.Cells.Item(135, 12).Value = 16794.6426
$ws.Cells.Item(135, 13).Value = -2063.0001
$ws.Cells.Item(135, 14).Value = -21864.6426

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4514.4443
$ws.Cells.Item(70, 9).Value = 4417.3335
$ws.Cells.Item(70, 10).Value = 5000
$ws.Cells.Item(70, 11).Value = 4417.3335
$ws.Cells.Item(70, 12).Value = 5000
$ws.Cells.Item(70, 13).Value = -4147.3335
$ws.Cells.Item(70, 14).Value = -5540
$ws.Cells.Item(73, 8).Value = 4514.4443
$ws.Cells.Item(73, 9).Value = 4417.3335
$ws.Cells.Item(73, 10).Value = 5000
$ws.Cells.Item(73, 11).Value = 4417.3335
$ws.Cells.Item(73, 12).Value = 5000
$ws.Cells.Item(73, 13).Value = -3481.3335
$ws.Cells.Item(73, 14).Value = -6872
$ws.Cells.Item(80, 8).Value = 4433.9736
$ws.Cells.Item(80, 9).Value = 5110.9614
$ws.Cells.Item(80, 10).Value = 2967.1667
$ws.Cells.Item(80, 11).Value = 5110.9614
$ws.Cells.Item(80, 12).Value = 2967.1667
$ws.Cells.Item(80, 13).Value = -4112.9614
$ws.Cells.Item(80, 14).Value = -4963.1667
$ws.Cells.Item(83, 8).Value = 4433.9736
$ws.Cells.Item(83, 9).Value = 5110.9614
$ws.Cells.Item(83, 10).Value = 2967.1667
$ws.Cells.Item(83, 11).Value = 25554.807
$ws.Cells.Item(83, 12).Value = 14835.8335
$ws.Cells.Item(83, 13).Value = -20562.807
$ws.Cells.Item(83, 14).Value = -24819.8335
$ws.Cells.Item(102, 8).Value = 3247.7917
$ws.Cells.Item(102, 9).Value = 2203.6155
$ws.Cells.Item(102, 10).Value = 4481.8184
$ws.Cells.Item(102, 11).Value = 2203.6155
$ws.Cells.Item(102, 12).Value = 4481.8184
$ws.Cells.Item(102, 13).Value = -581.6154999999999
$ws.Cells.Item(102, 14).Value = -7725.8184
$ws.Cells.Item(126, 8).Value = 2780.182
$ws.Cells.Item(126, 9).Value = 2748.2
$ws.Cells.Item(126, 10).Value = 3100
$ws.Cells.Item(126, 11).Value = 8244.599999999999
$ws.Cells.Item(126, 12).Value = 9300
$ws.Cells.Item(126, 13).Value = -5774.599999999999
$ws.Cells.Item(126, 14).Value = -14240
$ws.Cells.Item(132, 8).Value = 3409.1936
$ws.Cells.Item(132, 9).Value = 3499.2334
$ws.Cells.Item(132, 10).Value = 3324.7812
$ws.Cells.Item(132, 11).Value = 10497.7002
$ws.Cells.Item(132, 12).Value = 9974.3436
$ws.Cells.Item(132, 13).Value = -7967.700199999999
$ws.Cells.Item(132, 14).Value = -15034.3436

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 14163.5
$ws.Cells.Item(40, 9).Value = 15758.286
$ws.Cells.Item(40, 10).Value = 3000
$ws.Cells.Item(40, 11).Value = 15758.286
$ws.Cells.Item(40, 12).Value = 3000
$ws.Cells.Item(40, 13).Value = -15622.286
$ws.Cells.Item(40, 14).Value = -3272
$ws.Cells.Item(46, 8).Value = 1044.7727
$ws.Cells.Item(46, 10).Value = 614.5
$ws.Cells.Item(46, 12).Value = 614.5
$ws.Cells.Item(46, 14).Value = -990.5
$ws.Cells.Item(93, 8).Value = 711.25
$ws.Cells.Item(93, 9).Value = 658.2692
$ws.Cells.Item(93, 10).Value = 1400
$ws.Cells.Item(93, 11).Value = 658.2692
$ws.Cells.Item(93, 12).Value = 1400
$ws.Cells.Item(93, 13).Value = 589.7308
$ws.Cells.Item(93, 14).Value = -3896
$ws.Cells.Item(100, 8).Value = 76928856
$ws.Cells.Item(100, 9).Value = 8842.857
$ws.Cells.Item(100, 11).Value = 8842.857
$ws.Cells.Item(100, 13).Value = -8301.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2140.7576
$ws.Cells.Item(122, 9).Value = 2148.2693
$ws.Cells.Item(122, 10).Value = 2112.8572
$ws.Cells.Item(122, 11).Value = 6444.8079
$ws.Cells.Item(122, 12).Value = 6338.571599999999
$ws.Cells.Item(122, 13).Value = -3994.8079
$ws.Cells.Item(122, 14).Value = -11238.5716
$ws.Cells.Item(126, 8).Value = 1070.9615
$ws.Cells.Item(126, 9).Value = 812.25
$ws.Cells.Item(126, 10).Value = 1933.3334
$ws.Cells.Item(126, 11).Value = 2436.75
$ws.Cells.Item(126, 12).Value = 5800.0002
$ws.Cells.Item(126, 13).Value = 33.25
$ws.Cells.Item(126, 14).Value = -10740.0002
$ws.Cells.Item(132, 8).Value = 2384.9714
$ws.Cells.Item(132, 9).Value = 1320.375
$ws.Cells.Item(132, 10).Value = 3281.4736
$ws.Cells.Item(132, 11).Value = 3961.125
$ws.Cells.Item(132, 12).Value = 9844.4208
$ws.Cells.Item(132, 13).Value = -1431.125
$ws.Cells.Item(132, 14).Value = -14904.4208
